# Generate Report for Handoff
# Updates the localization-status report with the results of a new handoff:
# a freshly generated GUID-named markdown source file, new content-hash'd
# .xlf hand-off files, and refreshed handoff timestamps.

$wb = $excel.ActiveWorkbook

$oldGuid = "413d32dd-02e6-44c2-abd3-d05c68e017e4"
$newGuid = "9260208b-9501-4fc7-9428-44d9fe3aeb80"

$oldHash = "1f33a74e9e7961620735609a262682e2fdd4534f"
$newHash = "625d6ba34c428a04eaec17bcd5ce5d1b5fc03e1f"

$oldMdName = "$oldGuid.md"
$newMdName = "$newGuid.md"

$oldZhXlf = "$oldGuid.$oldHash.zh-cn.xlf"
$newZhXlf = "$newGuid.$newHash.zh-cn.xlf"

$oldDeXlf = "$oldGuid.$oldHash.de-de.xlf"
$newDeXlf = "$newGuid.$newHash.de-de.xlf"

# Hyperlink addresses stay pointing at their original targets - only the
# text shown in the cell / hyperlink display is refreshed.
$mdAddress = "https://github.com/OpenLocalizationTest/oltest/blob/27973f64ea8e0ef5c6fb3e43e4b5727ba7b016fb/e2e/$oldMdName"
$zhXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/62d656316b97b1e8fa50b3d233e09a38b789188a/ol-handoff/OpenLocalizationTestOrg/oltest.zh-cn/ci/ht/$oldZhXlf"
$deXlfAddress = "https://github.com/OpenLocalizationTestOrg/olhandoff/blob/e73af98a80b6c45104bf834415f441751ab35798/ol-handoff/OpenLocalizationTestOrg/oltest.de-de/ci/ht/$oldDeXlf"

$hyperlinkColor = 15570276 # BGR long for RGB FF6495ED (cornflower blue), matches the workbook's HyperLink style
$underlineSingle = 2       # xlUnderlineStyleSingle

function Set-LinkedCell($ws, $cellAddr, $linkAddress, $displayText) {
    $ws.Hyperlinks.Add($ws.Range($cellAddr), $linkAddress, "", "", $displayText)
    $ws.Range($cellAddr).Font.Underline = $underlineSingle
    $ws.Range($cellAddr).Font.Color = $hyperlinkColor
}

# ---------------------------------------------------------------------
# Sheet "Overview": File Name | zh-cn | de-de | Latest Handoff Date
# ---------------------------------------------------------------------
$wsOverview = $wb.Worksheets.Item("Overview")

$wsOverview.Hyperlinks.Delete()
Set-LinkedCell $wsOverview "A2" $mdAddress $newMdName

$wsOverview.Range("D2").Value = "2016-45-12 04:45:53"

# ---------------------------------------------------------------------
# Sheet "zh-cn"
# ---------------------------------------------------------------------
$wsZhCn = $wb.Worksheets.Item("zh-cn")

$wsZhCn.Hyperlinks.Delete()
Set-LinkedCell $wsZhCn "A2" $mdAddress $newMdName
Set-LinkedCell $wsZhCn "B2" $mdAddress ".md"
Set-LinkedCell $wsZhCn "D2" $zhXlfAddress $newZhXlf

$wsZhCn.Range("E2").Value = "2016-03-12 04:45:51"

# ---------------------------------------------------------------------
# Sheet "de-de"
# ---------------------------------------------------------------------
$wsDeDe = $wb.Worksheets.Item("de-de")

$wsDeDe.Hyperlinks.Delete()
Set-LinkedCell $wsDeDe "A2" $mdAddress $newMdName
Set-LinkedCell $wsDeDe "B2" $mdAddress ".md"
Set-LinkedCell $wsDeDe "D2" $deXlfAddress $newDeXlf

$wsDeDe.Range("E2").Value = "2016-03-12 04:45:53"
